# Scheduled price/profit refresh across the eight crafting-job sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# rows whose underlying market data moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1572.1
$ws.Range("I15").Value = 1572.1
$ws.Range("K15").Value = 4716.299999999999
$ws.Range("M15").Value = -4547.299999999999

$ws.Range("H40").Value = 2183.0857
$ws.Range("I40").Value = 2117.4546
$ws.Range("J40").Value = 2294.1538
$ws.Range("K40").Value = 2117.4546
$ws.Range("L40").Value = 2294.1538
$ws.Range("M40").Value = -1942.4546
$ws.Range("N40").Value = -2644.1538

$ws.Range("H41").Value = 833.3333
$ws.Range("I41").Value = 700
$ws.Range("K41").Value = 700
$ws.Range("M41").Value = -260

$ws.Range("H98").Value = 4548.4
$ws.Range("I98").Value = 1097.8
$ws.Range("K98").Value = 1097.8
$ws.Range("M98").Value = 400.2

$ws.Range("H100").Value = 1410.7273
$ws.Range("I100").Value = 1352.5714
$ws.Range("J100").Value = 1512.5
$ws.Range("K100").Value = 1352.5714
$ws.Range("L100").Value = 1512.5
$ws.Range("M100").Value = -811.5714
$ws.Range("N100").Value = -2594.5

$ws.Range("H106").Value = 2091
$ws.Range("I106").Value = 2091
$ws.Range("K106").Value = 2091
$ws.Range("M106").Value = -1460

$ws.Range("H112").Value = 806.9091
$ws.Range("J112").Value = 806.9091
$ws.Range("L112").Value = 2420.7273
$ws.Range("N112").Value = -4636.7273

$ws.Range("H122").Value = 4548.4
$ws.Range("I122").Value = 1097.8
$ws.Range("K122").Value = 3293.4
$ws.Range("M122").Value = -843.3999999999996

$ws.Range("H138").Value = 3380.1528
$ws.Range("J138").Value = 3410.3845
$ws.Range("L138").Value = 10231.1535
$ws.Range("N138").Value = -20511.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3151.8965
$ws.Range("I61").Value = 1295.85
$ws.Range("K61").Value = 1295.85
$ws.Range("M61").Value = -1083.85

$ws.Range("H74").Value = 2597.2666
$ws.Range("I74").Value = 2320
$ws.Range("J74").Value = 4399.5
$ws.Range("K74").Value = 2320
$ws.Range("L74").Value = 4399.5
$ws.Range("M74").Value = -1446
$ws.Range("N74").Value = -6147.5

$ws.Range("H77").Value = 2597.2666
$ws.Range("I77").Value = 2320
$ws.Range("J77").Value = 4399.5
$ws.Range("K77").Value = 11600
$ws.Range("L77").Value = 21997.5
$ws.Range("M77").Value = -7232
$ws.Range("N77").Value = -30733.5

$ws.Range("H88").Value = 1545.8462
$ws.Range("I88").Value = 866.3333
$ws.Range("K88").Value = 866.3333
$ws.Range("M88").Value = -460.3333

$ws.Range("H91").Value = 1545.8462
$ws.Range("I91").Value = 866.3333
$ws.Range("K91").Value = 866.3333
$ws.Range("M91").Value = 537.6667

$ws.Range("H136").Value = 3151.8965
$ws.Range("I136").Value = 1295.85
$ws.Range("K136").Value = 3887.55
$ws.Range("M136").Value = -1337.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2666.7856
$ws.Range("I86").Value = 873.5
$ws.Range("J86").Value = 3384.1
$ws.Range("K86").Value = 873.5
$ws.Range("L86").Value = 3384.1
$ws.Range("M86").Value = 249.5
$ws.Range("N86").Value = -5630.1

$ws.Range("H89").Value = 2666.7856
$ws.Range("I89").Value = 873.5
$ws.Range("J89").Value = 3384.1
$ws.Range("K89").Value = 4367.5
$ws.Range("L89").Value = 16920.5
$ws.Range("M89").Value = 1248.5
$ws.Range("N89").Value = -28152.5

$ws.Range("H94").Value = 612.8570999999999
$ws.Range("I94").Value = 615
$ws.Range("K94").Value = 615
$ws.Range("M94").Value = -164

$ws.Range("H96").Value = 9896.143
$ws.Range("I96").Value = 9896.143
$ws.Range("K96").Value = 9896.143
$ws.Range("M96").Value = -7150.143

$ws.Range("H99").Value = 3358.75
$ws.Range("I99").Value = 1949.1666
$ws.Range("J99").Value = 4768.3335
$ws.Range("K99").Value = 1949.1666
$ws.Range("L99").Value = 4768.3335
$ws.Range("M99").Value = -451.1666
$ws.Range("N99").Value = -7764.3335

$ws.Range("H107").Value = 4524.3076
$ws.Range("I107").Value = 4163.909
$ws.Range("J107").Value = 6506.5
$ws.Range("K107").Value = 4163.909
$ws.Range("L107").Value = 6506.5
$ws.Range("M107").Value = -2243.909
$ws.Range("N107").Value = -10346.5

$ws.Range("H122").Value = 79999
$ws.Range("J122").Value = 79999
$ws.Range("L122").Value = 79999
$ws.Range("N122").Value = -89799

$ws.Range("H134").Value = 792
$ws.Range("I134").Value = 794.4
$ws.Range("J134").Value = 780
$ws.Range("K134").Value = 2383.2
$ws.Range("L134").Value = 2340
$ws.Range("M134").Value = 151.8000000000002
$ws.Range("N134").Value = -7410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5603
$ws.Range("J31").Value = 5958.375
$ws.Range("L31").Value = 5958.375
$ws.Range("N31").Value = -6548.375

$ws.Range("H34").Value = 5603
$ws.Range("J34").Value = 5958.375
$ws.Range("L34").Value = 5958.375
$ws.Range("N34").Value = -6362.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18255076
$ws.Range("I4").Value = 20001330
$ws.Range("K4").Value = 60003990
$ws.Range("M4").Value = -60003878

$ws.Range("H5").Value = 416.66666
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0

$ws.Range("H23").Value = 15000396
$ws.Range("J23").Value = 475.8
$ws.Range("L23").Value = 1427.4
$ws.Range("N23").Value = -1897.4

$ws.Range("H56").Value = 9379.833000000001
$ws.Range("I56").Value = 9379.833000000001
$ws.Range("K56").Value = 9379.833000000001
$ws.Range("M56").Value = -8849.833000000001

$ws.Range("H121").Value = 14395.077
$ws.Range("J121").Value = 7092
$ws.Range("L121").Value = 21276
$ws.Range("N121").Value = -23896

$ws.Range("H122").Value = 608
$ws.Range("J122").Value = 567
$ws.Range("L122").Value = 5103
$ws.Range("N122").Value = -10003

$ws.Range("H134").Value = 18713.715
$ws.Range("I134").Value = 19000
$ws.Range("J134").Value = 18666
$ws.Range("K134").Value = 57000
$ws.Range("L134").Value = 55998
$ws.Range("M134").Value = -51930
$ws.Range("N134").Value = -66138

$ws.Range("H135").Value = 416.66666
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0

$ws.Range("H137").Value = 4928.2
$ws.Range("I137").Value = 4465
$ws.Range("K137").Value = 13395
$ws.Range("M137").Value = -8295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5254
$ws.Range("I70").Value = 4997.6665
$ws.Range("J70").Value = 6023
$ws.Range("K70").Value = 4997.6665
$ws.Range("L70").Value = 6023
$ws.Range("M70").Value = -4727.6665
$ws.Range("N70").Value = -6563

$ws.Range("H73").Value = 5254
$ws.Range("I73").Value = 4997.6665
$ws.Range("J73").Value = 6023
$ws.Range("K73").Value = 4997.6665
$ws.Range("L73").Value = 6023
$ws.Range("M73").Value = -4061.6665
$ws.Range("N73").Value = -7895

$ws.Range("H107").Value = 1353.2941
$ws.Range("I107").Value = 333.8
$ws.Range("J107").Value = 8999.5
$ws.Range("K107").Value = 333.8
$ws.Range("L107").Value = 8999.5
$ws.Range("M107").Value = 1586.2
$ws.Range("N107").Value = -12839.5

$ws.Range("H132").Value = 1997.6364
$ws.Range("I132").Value = 1664.1111
$ws.Range("K132").Value = 4992.3333
$ws.Range("M132").Value = -2462.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6860.5713
$ws.Range("I7").Value = 3807.111
$ws.Range("J7").Value = 8306.947
$ws.Range("K7").Value = 3807.111
$ws.Range("L7").Value = 8306.947
$ws.Range("M7").Value = -3695.111
$ws.Range("N7").Value = -8530.947

$ws.Range("H46").Value = 2618.3928
$ws.Range("I46").Value = 2034.8889
$ws.Range("J46").Value = 2894.7896
$ws.Range("K46").Value = 2034.8889
$ws.Range("L46").Value = 2894.7896
$ws.Range("M46").Value = -1846.8889
$ws.Range("N46").Value = -3270.7896

$ws.Range("H61").Value = 4787.3145
$ws.Range("I61").Value = 4292.6553
$ws.Range("K61").Value = 4292.6553
$ws.Range("M61").Value = -4090.6553

$ws.Range("H93").Value = 668.8333
$ws.Range("J93").Value = 864
$ws.Range("L93").Value = 864
$ws.Range("N93").Value = -3360

$ws.Range("H113").Value = 4787.3145
$ws.Range("I113").Value = 4292.6553
$ws.Range("K113").Value = 4292.6553
$ws.Range("M113").Value = -2122.6553

$ws.Range("H126").Value = 6860.5713
$ws.Range("I126").Value = 3807.111
$ws.Range("J126").Value = 8306.947
$ws.Range("K126").Value = 11421.333
$ws.Range("L126").Value = 24920.841
$ws.Range("M126").Value = -8951.332999999999
$ws.Range("N126").Value = -29860.841

$ws.Range("H132").Value = 3575.3076
$ws.Range("I132").Value = 2608.4
$ws.Range("J132").Value = 4179.625
$ws.Range("K132").Value = 7825.200000000001
$ws.Range("L132").Value = 12538.875
$ws.Range("M132").Value = -5295.200000000001
$ws.Range("N132").Value = -17598.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 52000
$ws.Range("J92").Value = 52000
$ws.Range("L92").Value = 52000
$ws.Range("N92").Value = -56992

$ws.Range("H122").Value = 1619
$ws.Range("I122").Value = 1619
$ws.Range("K122").Value = 4857
$ws.Range("M122").Value = -2407

$ws.Range("H132").Value = 3293.6667
$ws.Range("I132").Value = 3092.5715
$ws.Range("J132").Value = 3997.5
$ws.Range("K132").Value = 9277.7145
$ws.Range("L132").Value = 11992.5
$ws.Range("M132").Value = -6747.7145
$ws.Range("N132").Value = -17052.5

# Deletions (cells removed entirely)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N5").ClearContents()
$ws.Range("N135").ClearContents()
